# Form the consolidated report: fill in the "Absent" (column H) values
# for the rows that were left blank/incorrect, so that Absent = 1 - Real.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
